{"js": "// Replace the 100 arithmetic-problem strings in the single 20x5 table with\n// their new values, in document (row-major) order. The table's cell count\n// and grid shape (20 rows x 5 cols) do not change between before/after; only\n// the text content of each cell's run changes (one row's worth of cells was\n// dropped from the middle and the same number appended at the end, so every\n// cell position 1..100 maps to exactly one new value).\nconst NEW_VALUES = [\"63-27=\",\"58+26=\",\"49+5=\",\"38-29=\",\"84-7=\",\"82-19=\",\"76-8=\",\"46-7=\",\"71-3=\",\"11-3=\",\"60-8=\",\"50-2=\",\"59+28=\",\"90-28=\",\"64-6=\",\"68-19=\",\"69+6=\",\"26+56=\",\"21-7=\",\"94-76=\",\"33+38=\",\"41-28=\",\"76-39=\",\"76+6=\",\"81-53=\",\"60-5=\",\"59+8=\",\"24-7=\",\"57-9=\",\"53-46=\",\"58+15=\",\"19+47=\",\"49+29=\",\"93-78=\",\"72-5=\",\"37+47=\",\"53+29=\",\"40-17=\",\"43-7=\",\"58+36=\",\"74-67=\",\"73-29=\",\"51-14=\",\"5+69=\",\"70-2=\",\"26+49=\",\"56+38=\",\"96-68=\",\"39+16=\",\"8+88=\",\"17+18=\",\"25+26=\",\"37+27=\",\"48+48=\",\"74-9=\",\"17+46=\",\"59+17=\",\"3+48=\",\"56+35=\",\"7+29=\",\"18+54=\",\"29+34=\",\"65+6=\",\"5+36=\",\"54-35=\",\"22-4=\",\"76+9=\",\"91-13=\",\"76-7=\",\"58+19=\",\"63-45=\",\"4+18=\",\"5+37=\",\"15+17=\",\"66-17=\",\"32-26=\",\"65+6=\",\"82-64=\",\"94-59=\",\"86-7=\",\"69+25=\",\"46+28=\",\"26+47=\",\"29+26=\",\"71-38=\",\"97-18=\",\"91-43=\",\"35+49=\",\"50-29=\",\"50-8=\",\"31-13=\",\"48+18=\",\"8+39=\",\"35+19=\",\"86-18=\",\"87-39=\",\"88-69=\",\"15+6=\",\"62-58=\",\"38+9=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's current value up front (one sync) so we only write the\n// ones that actually changed.\nconst cellGrid = [];\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  cellGrid.push(cells);\n}\nawait context.sync();\n\nfor (let r = 0; r < cellGrid.length; r++) {\n  cellGrid[r].items.forEach((cell) => cell.load(\"value\"));\n}\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < cellGrid.length; r++) {\n  const cellItems = cellGrid[r].items;\n  for (let c = 0; c < cellItems.length; c++) {\n    const newVal = NEW_VALUES[idx];\n    if (newVal !== undefined && cellItems[c].value !== newVal) {\n      cellItems[c].value = newVal;\n    }\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem strings in the single 20x5 table with\n# their new values, in document (row-major) order. The table's cell count\n# and grid shape (20 rows x 5 cols) do not change between before/after; only\n# the text content of each cell's run changes (one row's worth of cells was\n# dropped from the middle and the same number appended at the end, so every\n# cell position 1..100 maps to exactly one new value).\n$newValues = @(\"63-27=\", \"58+26=\", \"49+5=\", \"38-29=\", \"84-7=\", \"82-19=\", \"76-8=\", \"46-7=\", \"71-3=\", \"11-3=\", \"60-8=\", \"50-2=\", \"59+28=\", \"90-28=\", \"64-6=\", \"68-19=\", \"69+6=\", \"26+56=\", \"21-7=\", \"94-76=\", \"33+38=\", \"41-28=\", \"76-39=\", \"76+6=\", \"81-53=\", \"60-5=\", \"59+8=\", \"24-7=\", \"57-9=\", \"53-46=\", \"58+15=\", \"19+47=\", \"49+29=\", \"93-78=\", \"72-5=\", \"37+47=\", \"53+29=\", \"40-17=\", \"43-7=\", \"58+36=\", \"74-67=\", \"73-29=\", \"51-14=\", \"5+69=\", \"70-2=\", \"26+49=\", \"56+38=\", \"96-68=\", \"39+16=\", \"8+88=\", \"17+18=\", \"25+26=\", \"37+27=\", \"48+48=\", \"74-9=\", \"17+46=\", \"59+17=\", \"3+48=\", \"56+35=\", \"7+29=\", \"18+54=\", \"29+34=\", \"65+6=\", \"5+36=\", \"54-35=\", \"22-4=\", \"76+9=\", \"91-13=\", \"76-7=\", \"58+19=\", \"63-45=\", \"4+18=\", \"5+37=\", \"15+17=\", \"66-17=\", \"32-26=\", \"65+6=\", \"82-64=\", \"94-59=\", \"86-7=\", \"69+25=\", \"46+28=\", \"26+47=\", \"29+26=\", \"71-38=\", \"97-18=\", \"91-43=\", \"35+49=\", \"50-29=\", \"50-8=\", \"31-13=\", \"48+18=\", \"8+39=\", \"35+19=\", \"86-18=\", \"87-39=\", \"88-69=\", \"15+6=\", \"62-58=\", \"38+9=\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $newVal = $newValues[$idx]\n        $cell = $t.Cell($r, $c)\n        if ($cell.Range.Text -ne $newVal) {\n            $cell.Range.Text = $newVal\n        }\n        $idx++\n    }\n}\n"}
